# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update per-country COVID figures (which also changes the sort order
#   for a few countries whose totals overtook their neighbours)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 12:42"

# --- Row 18/19: Banglades overtakes Pakistan ---
$ws.Range("A18").Value = "Banglades"
$ws.Range("B18").Value = 292625
$ws.Range("C18").Value = 2265
$ws.Range("D18").Value = 175567
$ws.Range("E18").Value = 113151
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 46
$ws.Range("H18").Value = 3907

$ws.Range("A19").Value = "Pakistan"
$ws.Range("B19").Value = 292174
$ws.Range("C19").Value = 586
$ws.Range("D19").Value = 275317
$ws.Range("E19").Value = 10626
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = 6231

# --- Row 42: Rumania figures updated (no rank change) ---
$ws.Range("A42").Value = "Rumania"
$ws.Range("B42").Value = 77544
$ws.Range("C42").Value = 1189
$ws.Range("D42").Value = 35079
$ws.Range("E42").Value = 39232
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 37
$ws.Range("H42").Value = 3233

# --- Rows 121-123: Eslovaquia overtakes Somalia and Mayotte ---
$ws.Range("A121").Value = "Eslovaquia"
$ws.Range("B121").Value = 3316
$ws.Range("C121").Value = 91
$ws.Range("D121").Value = 2147
$ws.Range("E121").Value = 1136
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 33

$ws.Range("A122").Value = "Somalia"
$ws.Range("B122").Value = 3265
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 2396
$ws.Range("E122").Value = 776
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 93

$ws.Range("A123").Value = "Mayotte"
$ws.Range("B123").Value = 3237
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 2964
$ws.Range("E123").Value = 234
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 39

# --- Row 133: Estonia figures updated (no rank change) ---
$ws.Range("A133").Value = "Estonia"
$ws.Range("B133").Value = 2265
$ws.Range("C133").Value = 21
$ws.Range("D133").Value = 2024
$ws.Range("E133").Value = 178
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 63

# --- Row 208: Dominica figures updated (no rank change) ---
$ws.Range("A208").Value = "Dominica"
$ws.Range("B208").Value = 19
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 18
$ws.Range("E208").Value = 1
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
